$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Remove the trailing rows (old rows 8-15) so the sheet shrinks to A1:H7 ---
$ws.Range("A8:H15").EntireRow.Delete()

# --- Column width tweaks (B: 55->51, D: 32->28, H: 16->12) ---
# Excel's ColumnWidth COM property stores ~0.8333 (5/6) wider than the raw
# OOXML <col width> value, so back that padding out to land on the exact
# target stored width.
$ws.Columns.Item(2).ColumnWidth = 51 - 5/6
$ws.Columns.Item(4).ColumnWidth = 28 - 5/6
$ws.Columns.Item(8).ColumnWidth = 12 - 5/6

# --- Row 2 ---
$ws.Range("A2").Value = "2025-11-19 06:26:44"
$ws.Range("B2").Value = "【せどり×ツール製作】APIを使用したせどりツールを製作できるエンジニアさんを募集します♪"
$ws.Range("C2").Value = "システム開発"
$ws.Range("D2").Value = "20,000 円 ~ 50,000 円 / 固定"
$ws.Range("E2").Value = "期限情報なし"
$ws.Range("F2").Value = "https://www.lancers.jp/work/detail/5217096"
$ws.Range("G2").Value = 243
$ws.Range("H2").Value = "🔥API ◆ツール"

# --- Row 3 ---
$ws.Range("A3").Value = "2025-11-19 06:26:44"
$ws.Range("B3").Value = "ホットペッパービューティーブログ一括投稿システム開発"
$ws.Range("C3").Value = "システム開発"
$ws.Range("D3").Value = "20,000 円 ~ 50,000 円 / 固定"
$ws.Range("E3").Value = "期限情報なし"
$ws.Range("F3").Value = "https://www.lancers.jp/work/detail/5437096"
$ws.Range("G3").Value = 113
$ws.Range("H3").Value = "◆開発,システム開発"

# --- Row 4 ---
$ws.Range("A4").Value = "2025-11-19 06:26:44"
$ws.Range("B4").Value = "【急募】Android用のライブ壁紙アプリ開発エンジニアを探しています!"
$ws.Range("C4").Value = "システム開発"
$ws.Range("D4").Value = "50,000 円 ~ 100,000 円 / 固定"
$ws.Range("E4").Value = "期限情報なし"
$ws.Range("F4").Value = "https://www.lancers.jp/work/detail/5436829"
$ws.Range("G4").Value = 93
$ws.Range("H4").Value = "◆開発 ◇アプリ"

# --- Row 5 ---
$ws.Range("A5").Value = "2025-11-19 06:26:44"
$ws.Range("B5").Value = "【急募】業務システムの要件定義と基本設計ができる方"
$ws.Range("C5").Value = "システム開発"
$ws.Range("D5").Value = "100,000 円 ~ 200,000 円 / 固定"
$ws.Range("E5").Value = "期限情報なし"
$ws.Range("F5").Value = "https://www.lancers.jp/work/detail/5437193"
$ws.Range("G5").Value = 33
$ws.Range("H5").ClearContents()

# --- Row 6 ---
$ws.Range("A6").Value = "2025-11-19 06:26:44"
$ws.Range("B6").Value = "【急募】行政向けシステム保守・運用の専門家を募集!"
$ws.Range("C6").Value = "システム開発"
$ws.Range("D6").Value = "200,000 円 ~ 300,000 円 / 固定"
$ws.Range("E6").Value = "期限情報なし"
$ws.Range("F6").Value = "https://www.lancers.jp/work/detail/5437146"
$ws.Range("G6").Value = 33
$ws.Range("H6").ClearContents()

# --- Row 7 ---
$ws.Range("A7").Value = "2025-11-19 06:26:44"
$ws.Range("B7").Value = "yahooプレイス用Worepressプラグイン 投稿記事をyahooプレイスのお知らせに要約投稿"
$ws.Range("C7").Value = "システム開発"
$ws.Range("D7").Value = "20,000 円 ~ 50,000 円 / 固定"
$ws.Range("E7").Value = "期限情報なし"
$ws.Range("F7").Value = "https://www.lancers.jp/work/detail/5436950"
$ws.Range("G7").Value = 13
$ws.Range("H7").ClearContents()

# --- Rebuild hyperlinks on column F from scratch (row deletion above leaves ---
# --- stale hyperlink registrations behind, so wipe then re-add cleanly)     ---
$ws.Cells.Hyperlinks.Delete()
$ws.Hyperlinks.Add($ws.Range("F2"), "https://www.lancers.jp/work/detail/5217096")
$ws.Hyperlinks.Add($ws.Range("F3"), "https://www.lancers.jp/work/detail/5437096")
$ws.Hyperlinks.Add($ws.Range("F4"), "https://www.lancers.jp/work/detail/5436829")
$ws.Hyperlinks.Add($ws.Range("F5"), "https://www.lancers.jp/work/detail/5437193")
$ws.Hyperlinks.Add($ws.Range("F6"), "https://www.lancers.jp/work/detail/5437146")
$ws.Hyperlinks.Add($ws.Range("F7"), "https://www.lancers.jp/work/detail/5436950")
